$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally has 9 data rows (rows 2-10): the sending/target cluster
# combination covered ECs/FAPs/MuSCs x ECs/FAPs/MuSCs (minus a couple). The
# updated data (recomputed from new TPM values) only has 6 data rows, covering
# sending clusters ECs/FAPs/MuSCs against target clusters FAPs/MuSCs only.
# Remove the last three rows (old rows 8,9,10) so the sheet shrinks from
# A1:T10 down to A1:T7.
$ws.Rows.Item(10).Delete() | Out-Null
$ws.Rows.Item(9).Delete() | Out-Null
$ws.Rows.Item(8).Delete() | Out-Null

# New data for rows 2-7, columns A-T (updated TPM-derived values).
$data = @(
    @("ECs","Icosl","Icos","FAPs",3,1,0.5407596666666666,1.622279,0.03618231591230665,0.03618231591230665,1,0.3333333333333333,0.107019,0.321057,0.2365534751714524,0.2365534751714524,0.05787155876699999,0.5208440289029999,0.00855905256880748,0.00855905256880748),
    @("ECs","Icosl","Icos","MuSCs",3,1,0.5407596666666666,1.622279,0.03618231591230665,0.03618231591230665,3,1,0.3453903333333333,1.036171,0.7634465248285476,0.7634465248285476,0.1867731615232222,1.680958453709,0.02762326334349917,0.02762326334349917),
    @("FAPs","Icosl","Icos","FAPs",3,1,11.32416633333333,33.972499,0.7577017832003755,0.7577017832003754,1,0.3333333333333333,0.107019,0.321057,0.2365534751714524,0.2365534751714524,1.211900956827,10.907108611443,0.1792369899596553,0.1792369899596552),
    @("FAPs","Icosl","Icos","MuSCs",3,1,11.32416633333333,33.972499,0.7577017832003755,0.7577017832003754,3,1,0.3453903333333333,1.036171,0.7634465248285476,0.7634465248285476,3.91125758459211,35.201318261329,0.5784647932407203,0.5784647932407202),
    @("MuSCs","Icosl","Icos","FAPs",3,1,3.080487333333333,9.241461999999999,0.2061159008873179,0.2061159008873179,1,0.3333333333333333,0.107019,0.321057,0.2365534751714524,0.2365534751714524,0.3296706739259999,2.967036065333999,0.04875743264298971,0.04875743264298971),
    @("MuSCs","Icosl","Icos","MuSCs",3,1,3.080487333333333,9.241461999999999,0.2061159008873179,0.2061159008873179,3,1,0.3453903333333333,1.036171,0.7634465248285476,0.7634465248285476,1.063970546889111,9.575734922001997,0.1573584682443282,0.1573584682443282)
)

$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
